$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.414.56'
$ws.Range('E2').Value = '  +0.77%  '

$ws.Range('D3').Value = '1.868.14'
$ws.Range('E3').Value = '  +0.34%  '

$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('D5').Value = '246.53'
$ws.Range('E5').Value = '  +1.96%  '

$ws.Range('E6').Value = '  +0.14%  '

$ws.Range('D7').Value = '0.4731'
$ws.Range('E7').Value = '  +0.25%  '

$ws.Range('D8').Value = '0.2902'
$ws.Range('E8').Value = '  +1.70%  '

$ws.Range('D9').Value = '0.06486'
$ws.Range('E9').Value = '  +0.20%  '

$ws.Range('D10').Value = '21.93'
$ws.Range('E10').Value = '  +5.86%  '

$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = '0.07711'
$ws.Range('E11').Value = '  +0.53%  '

$ws.Range('B12').Value = 'Litecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D12').Value = '97.68'
$ws.Range('E12').Value = '  +3.97%  '

$ws.Range('D13').Value = '0.7349'
$ws.Range('E13').Value = '  +8.02%  '

$ws.Range('D14').Value = '1.871.39'
$ws.Range('E14').Value = '  +0.52%  '

$ws.Range('D15').Value = '5.101'
$ws.Range('E15').Value = '  +0.76%  '

$ws.Range('D16').Value = '273.14'
$ws.Range('E16').Value = '  +1.62%  '

$ws.Range('D17').Value = '30.396.19'
$ws.Range('E17').Value = '  +0.72%  '

$ws.Range('D18').Value = '13.36'
$ws.Range('E18').Value = '  +0.12%  '

$ws.Range('D19').Value = '0.000007532'
$ws.Range('E19').Value = '  -0.16%  '

$ws.Range('E20').Value = '  +0.07%  '

$ws.Range('D21').Value = '2.117.24'
$ws.Range('E21').Value = '  +0.70%  '

$ws.Range('E22').Value = '  +0.20%  '

$ws.Range('D23').Value = '5.216'
$ws.Range('E23').Value = '  +0.91%  '

$ws.Range('D24').Value = '6.154'
$ws.Range('E24').Value = '  +0.97%  '

$ws.Range('D25').Value = '9.251'
$ws.Range('E25').Value = '  -0.80%  '

$ws.Range('D26').Value = '163.69'
$ws.Range('E26').Value = '  -1.40%  '

$ws.Range('D27').Value = '18.79'
$ws.Range('E27').Value = '  +0.42%  '

$ws.Range('D28').Value = '1.923'
$ws.Range('E28').Value = '  +2.18%  '

$ws.Range('D29').Value = '0.1002'
$ws.Range('E29').Value = '  +1.64%  '

$ws.Range('E30').Value = '  -0.78%  '

$ws.Range('D31').Value = '1.507'
$ws.Range('E31').Value = '  +0.14%  '

$ws.Range('D32').Value = '4.294'
$ws.Range('E32').Value = '  +1.63%  '

$ws.Range('D33').Value = '4.131'
$ws.Range('E33').Value = '  +3.40%  '

$ws.Range('D34').Value = '0.04809'
$ws.Range('E34').Value = '  +2.41%  '

$ws.Range('D35').Value = '1.116'
$ws.Range('E35').Value = '  +0.62%  '

$ws.Range('D36').Value = '0.6958'
$ws.Range('E36').Value = '  +1.53%  '

$ws.Range('D37').Value = '2.712'
$ws.Range('E37').Value = '  +0.17%  '

$ws.Range('D38').Value = '0.01852'
$ws.Range('E38').Value = '  +1.26%  '

$ws.Range('D39').Value = '2.748'
$ws.Range('E39').Value = '  +1.02%  '

$ws.Range('D40').Value = '6.302'
$ws.Range('E40').Value = '  -1.01%  '

$ws.Range('D41').Value = '1.964'
$ws.Range('E41').Value = '  +4.06%  '

$ws.Range('D42').Value = '71.41'
$ws.Range('E42').Value = '  +1.95%  '

$ws.Range('D43').Value = '0.4175'
$ws.Range('E43').Value = '  +2.84%  '

$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.16%  '

$ws.Range('D45').Value = '0.8339'
$ws.Range('E45').Value = '  -0.25%  '

$ws.Range('D46').Value = '102.50'
$ws.Range('E46').Value = '  +0.51%  '

$ws.Range('D47').Value = '9.246'
$ws.Range('E47').Value = '  +0.21%  '

$ws.Range('D48').Value = '6.995'
$ws.Range('E48').Value = '  +1.06%  '

$ws.Range('D49').Value = '35.31'
$ws.Range('E49').Value = '  +2.96%  '

$ws.Range('D50').Value = '916.91'
$ws.Range('E50').Value = '  -1.03%  '

$ws.Range('D51').Value = '0.05632'
$ws.Range('E51').Value = '  +1.39%  '

